$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C1').Value = 'orgRole'
$ws.Range('B10').Value = 'For a user belonging to an organization with a Preliminary Access role, hide all verySensitive class information, as though it does not exist, regardless of the medium or location.'
$ws.Range('B11').Value = 'For a user belonging to an organization with a Preliminary Access role, hide all orgSensitive class information, as though it does not exist, regardless of the medium or location.'
$ws.Range('B12').Value = 'For a user belonging to an organization with a Preliminary Access role, hide all sensitive class information, as though it does not exist, regardless of the medium or location.'
$ws.Range('B13').Value = 'For a user belonging to an organization with a Preliminary Access role, truncate all zip codes to five digits, and replace numbers with Xs in all other fields (e.g. phone numbers and address numbers), regardless of the medium or location.'
$ws.Range('B15').Value = 'For a user belonging to an organization with a Preliminary Access role, display all public class information, regardless of the medium or location.'
$ws.Range('B16').Value = 'For a user belonging to an organization with a Preliminary Access role, hide all information by default, if it does not fall into one of the previous rules, regardless of the medium or location.'
$ws.Range('B17').Value = 'For a user belonging to an organization with a Statistics Access role, hide all verySensitive class information, as though it does not exist, regardless of the medium or location.'
$ws.Range('B18').Value = 'For a user belonging to an organization with a Statistics Access role, hide all orgSensitive class information, as though it does not exist, regardless of the medium or location.'
$ws.Range('B19').Value = 'For a user belonging to an organization with a Statistics Access role, hide all sensitive class information, as though it does not exist, regardless of the medium or location.'
$ws.Range('B20').Value = 'For a user belonging to an organization with a Statistics Access role, truncate all zip codes to five digits, and replace numbers with Xs in all other lessSensitive class fields (e.g. phone numbers and address numbers), regardless of the medium or location.'
$ws.Range('B22').Value = 'For a user belonging to an organization with a Statistics Access role, display all public class information, regardless of the medium or location.'
$ws.Range('B23').Value = 'For a user belonging to an organization with a Statistics Access role, hide all information by default, if it does not fall into one of the previous rules, regardless of the medium or location.'
$ws.Range('B24').Value = 'For a user belonging to an organization with a Situational Awareness role, replace verySensitive class information with a translated message if the case is within the users geofence, regardless of the medium if the user has not claimed or reported the case. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B26').Value = 'For a user belonging to an organization with a Situational Awareness role, show all case data for any other sensitivity class if the case is inside the users geofence, regardless of the medium, if the user has not claimed or reported the case.'
$ws.Range('B27').Value = 'For a user belonging to an organization with a Situational Awareness role, replace verySensitive class information with a translated message if the case is outside the users geofence, regardless of the medium, if the user has not claimed or reported the case. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B29').Value = 'For a user belonging to an organization with a Situational Awareness role, convert orgSensitive class information to a boolean true/false if the case is outside the users geofence, regardless of the medium, if the user has not claimed or reported the case. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B31').Value = 'For a user belonging to an organization with a Situational Awareness role, hide all sensitive class information, as though it does not exist, if the case is outside the users geofence, regardless of the medium, if the user has not claimed or reported the case. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B33').Value = 'For a user belonging to an organization with a Situational Awareness role, truncate all zip codes to five digits, and replace numbers with Xs in all other lessSensitive class fields (e.g. phone numbers and address numbers), if the case is outside the users geofence and the user has not claimed or reported the case, regardless of the medium. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B36').Value = 'For a user belonging to an organization with a Situational Awareness role, show all public class information, if the case is outside the users geofence, regardless of the medium. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B37').Value = 'For a user belonging to an organization with a Situational Awareness role, hide all information by default, if it does not fall into one of the previous rules, if the case is outside the users geofence, regardless of the medium. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B39').Value = 'For a user belonging to an organization with a Coordination Access role, replace verySensitive class information with a translated message, regardless of the location or medium, if the user has not claimed or reported the case. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B41').Value = 'For a user belonging to an organization with a Coordination Access role, convert orgSensitive class information to a boolean true/false, regardless of the location or medium, if the user has not claimed or reported the case. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B43').Value = 'For a user belonging to an organization with a Coordination Access role, hide all sensitive class information, as though it does not exist, regardless of the location or medium, if the user has not claimed or reported the case. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B45').Value = 'For a user belonging to an organization with a Coordination Access role, truncate all zip codes to five digits, and replace numbers with Xs in all other lessSensitive class fields (e.g. phone numbers and address numbers), regardless of the location or medium. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B48').Value = 'For a user belonging to an organization with a Coordination Access role, show all public class information, regardless of the location or medium. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B49').Value = 'For a user belonging to an organization with a Situational Awareness role, hide all information by default, if it does not fall into one of the previous rules, regardless of the medium or location. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B51').Value = 'For a user belonging to an organization with an LTR Access role, show all case information inside the users geofence, regardless of who reported or claimed it, and regardless of the medium.'
$ws.Range('B52').Value = 'For a user belonging to an organization with an LTR Access role, replace verySensitive class information with a translated message, if the case is outside the users geocode, regardless of the medium, if the user has not claimed or reported the case. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B54').Value = 'For a user belonging to an organization with an LTR Access role, convert orgSensitive class information to a boolean true/false, if the case is outside the users geocode, regardless of the medium, if the user has not claimed or reported the case. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B56').Value = 'For a user belonging to an organization with an LTR Access role, hide all sensitive class information, as though it does not exist, if the case is outside the users geocode, regardless of the medium, if the user has not claimed or reported the case. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B58').Value = 'For a user belonging to an organization with an LTR Access role, truncate all zip codes to five digits, and replace numbers with Xs in all other lessSensitive class fields (e.g. phone numbers and address numbers), if the case is outside the users geocode, regardless of the medium. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B61').Value = 'For a user belonging to an organization with an LTR Access role, show all public class information, if the case is outside the users geocode, regardless of the medium. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B62').Value = 'For a user belonging to an organization with an LTR Access role, hide all information by default, if it does not fall into one of the previous rules, if the case is outside the users geocode, regardless of the medium. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B64').Value = 'For a user belonging to an organization with a Recovery Access role, replace verySensitive class information with a translated message, regardless of the location, if the medium is download, if the user has not claimed or reported the case. However, inherit any higher level access from affiliated organizations.'
$ws.Range('B66').Value = 'For a user belonging to an organization with a Recovery Access role, show all other case information, regardless of the location, if the medium is download, if the user has not claimed or reported the case.'
$ws.Range('B67').Value = 'For a user belonging to an organization with a Recovery Access role, show all other case information, regardless of the location, if the medium is anything other than download, if the user has not claimed or reported the case.'
